$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the formulas in J2 and J3 - they were inverted (D/E instead of E/D)
$ws.Range("J2").Formula = "=E2/D2*100"
$ws.Range("J3").Formula = "=E3/D3*100"

# Update the selected cell in the sheet view
$ws.Range("J4").Select()
